$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3158.111
$ws.Range("I40").Value = 2360.0454
$ws.Range("J40").Value = 4412.2144
$ws.Range("K40").Value = 2360.0454
$ws.Range("L40").Value = 4412.2144
$ws.Range("M40").Value = -2185.0454
$ws.Range("N40").Value = -4762.2144

$ws.Range("H58").Value = 1204.2222
$ws.Range("I58").Value = 106.666664
$ws.Range("K58").Value = 319.999992
$ws.Range("M58").Value = -169.999992

$ws.Range("H64").Value = 7676.846
$ws.Range("I64").Value = 5133.1665
$ws.Range("K64").Value = 5133.1665
$ws.Range("M64").Value = -4885.1665

$ws.Range("H67").Value = 7676.846
$ws.Range("I67").Value = 5133.1665
$ws.Range("K67").Value = 5133.1665
$ws.Range("M67").Value = -4275.1665

$ws.Range("H70").Value = 3221.1
$ws.Range("I70").Value = 2033
$ws.Range("J70").Value = 3730.2856
$ws.Range("K70").Value = 6099
$ws.Range("L70").Value = 11190.8568
$ws.Range("M70").Value = -5829
$ws.Range("N70").Value = -11730.8568

$ws.Range("H73").Value = 3221.1
$ws.Range("I73").Value = 2033
$ws.Range("J73").Value = 3730.2856
$ws.Range("K73").Value = 6099
$ws.Range("L73").Value = 11190.8568
$ws.Range("M73").Value = -5163
$ws.Range("N73").Value = -13062.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27784758
$ws.Range("I32").Value = 27784758
$ws.Range("K32").Value = 27784758
$ws.Range("M32").Value = -27784471

$ws.Range("H74").Value = 1822.5834
$ws.Range("I74").Value = 1816.5238
$ws.Range("K74").Value = 1816.5238
$ws.Range("M74").Value = -942.5237999999999

$ws.Range("H77").Value = 1822.5834
$ws.Range("I77").Value = 1816.5238
$ws.Range("K77").Value = 9082.618999999999
$ws.Range("M77").Value = -4714.618999999999

$ws.Range("H132").Value = 4355.9375
$ws.Range("J132").Value = 4779.1816
$ws.Range("L132").Value = 14337.5448
$ws.Range("N132").Value = -19397.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H134").Value = 4790.724
$ws.Range("I134").Value = 3791.5264
$ws.Range("K134").Value = 11374.5792
$ws.Range("M134").Value = -8839.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3783.8
$ws.Range("J31").Value = 4187.9473
$ws.Range("L31").Value = 4187.9473
$ws.Range("N31").Value = -4777.9473

$ws.Range("H34").Value = 3783.8
$ws.Range("J34").Value = 4187.9473
$ws.Range("L34").Value = 4187.9473
$ws.Range("N34").Value = -4591.9473

$ws.Range("H58").Value = 6279.0713
$ws.Range("I58").Value = 4308
$ws.Range("J58").Value = 7067.5
$ws.Range("K58").Value = 4308
$ws.Range("L58").Value = 7067.5
$ws.Range("M58").Value = -4105
$ws.Range("N58").Value = -7473.5

$ws.Range("H62").Value = 4109.1113
$ws.Range("I62").Value = 4392.75
$ws.Range("J62").Value = 1840
$ws.Range("K62").Value = 4392.75
$ws.Range("L62").Value = 1840
$ws.Range("M62").Value = -3768.75
$ws.Range("N62").Value = -3088

$ws.Range("H65").Value = 4109.1113
$ws.Range("I65").Value = 4392.75
$ws.Range("J65").Value = 1840
$ws.Range("K65").Value = 21963.75
$ws.Range("L65").Value = 9200
$ws.Range("M65").Value = -18843.75
$ws.Range("N65").Value = -15440

$ws.Range("H105").Value = 8970.182000000001
$ws.Range("I105").Value = 1279
$ws.Range("K105").Value = 1279
$ws.Range("M105").Value = 468

$ws.Range("H107").Value = 8519.308000000001
$ws.Range("I107").Value = 913.8
$ws.Range("K107").Value = 913.8
$ws.Range("M107").Value = 1006.2

$ws.Range("H119").Value = 42494.5
$ws.Range("J119").Value = 42494.5
$ws.Range("L119").Value = 42494.5
$ws.Range("N119").Value = -52170.5

$ws.Range("H121").Value = 19998.9
$ws.Range("J121").Value = 19998.9
$ws.Range("L121").Value = 19998.9
$ws.Range("N121").Value = -22618.9

$ws.Range("H132").Value = 3267.111
$ws.Range("I132").Value = 1917.5
$ws.Range("K132").Value = 5752.5
$ws.Range("M132").Value = -3222.5

$ws.Range("H134").Value = 5786.853
$ws.Range("I134").Value = 4209.7856
$ws.Range("K134").Value = 12629.3568
$ws.Range("M134").Value = -10094.3568

$ws.Range("H136").Value = 6279.0713
$ws.Range("I136").Value = 4308
$ws.Range("J136").Value = 7067.5
$ws.Range("K136").Value = 12924
$ws.Range("L136").Value = 21202.5
$ws.Range("M136").Value = -10374
$ws.Range("N136").Value = -26302.5

$ws.Range("H138").Value = 67856.78999999999
$ws.Range("J138").Value = 69230.766
$ws.Range("L138").Value = 69230.766
$ws.Range("N138").Value = -79510.766

$ws.Range("H139").Value = 47332.223
$ws.Range("J139").Value = 47332.223
$ws.Range("L139").Value = 47332.223
$ws.Range("N139").Value = -57612.223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 431.8
$ws.Range("J107").Value = 431.8
$ws.Range("L107").Value = 1295.4
$ws.Range("N107").Value = -5135.4

$ws.Range("H113").Value = 959.3889
$ws.Range("I113").Value = 779.7143
$ws.Range("J113").Value = 1073.7273
$ws.Range("K113").Value = 2339.1429
$ws.Range("L113").Value = 3221.1819
$ws.Range("M113").Value = -169.1428999999998
$ws.Range("N113").Value = -7561.1819

$ws.Range("H139").Value = 3091.1333
$ws.Range("I139").Value = 1613.9166
$ws.Range("J139").Value = 9000
$ws.Range("K139").Value = 4841.7498
$ws.Range("L139").Value = 27000
$ws.Range("M139").Value = 298.2502000000004
$ws.Range("N139").Value = -37280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15153796
$ws.Range("I80").Value = 2924.125
$ws.Range("K80").Value = 2924.125
$ws.Range("M80").Value = -1926.125

$ws.Range("H83").Value = 15153796
$ws.Range("I83").Value = 2924.125
$ws.Range("K83").Value = 14620.625
$ws.Range("M83").Value = -9628.625

$ws.Range("H132").Value = 2999.75
$ws.Range("I132").Value = 2999.75
$ws.Range("K132").Value = 8999.25
$ws.Range("M132").Value = -6469.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 57499
$ws.Range("J36").Value = 57499
$ws.Range("L36").Value = 57499
$ws.Range("N36").Value = -58623

$ws.Range("H40").Value = 7926.963
$ws.Range("I40").Value = 9769.154
$ws.Range("K40").Value = 9769.154
$ws.Range("M40").Value = -9633.154

$ws.Range("H55").Value = 1525.375
$ws.Range("I55").Value = 2409.4443
$ws.Range("K55").Value = 2409.4443
$ws.Range("M55").Value = -2236.4443

$ws.Range("H131").Value = 69929
$ws.Range("J131").Value = 69929
$ws.Range("L131").Value = 69929
$ws.Range("N131").Value = -80009

$ws.Range("H140").Value = 76000
$ws.Range("J140").Value = 76000
$ws.Range("L140").Value = 76000
$ws.Range("N140").Value = -86360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6946499
$ws.Range("I81").Value = 2295.0908
$ws.Range("J81").Value = 22223746
$ws.Range("K81").Value = 4590.1816
$ws.Range("L81").Value = 44447492
$ws.Range("M81").Value = -3529.1816
$ws.Range("N81").Value = -44449614

$ws.Range("H84").Value = 6946499
$ws.Range("I84").Value = 2295.0908
$ws.Range("J84").Value = 22223746
$ws.Range("K84").Value = 22950.908
$ws.Range("L84").Value = 222237460
$ws.Range("M84").Value = -17646.908
$ws.Range("N84").Value = -222248068

$ws.Range("H132").Value = 2832.7896
$ws.Range("I132").Value = 1544.6666
$ws.Range("K132").Value = 4633.9998
$ws.Range("M132").Value = -2103.9998

$ws.Range("H138").Value = 84075
$ws.Range("J138").Value = 84075
$ws.Range("L138").Value = 84075
$ws.Range("N138").Value = -94355

$ws.Range("H141").Value = 85416.664
$ws.Range("J141").Value = 85416.664
$ws.Range("L141").Value = 85416.664
$ws.Range("N141").Value = -95776.664
